$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Card Back" image filenames (rows 200-212, column E).
#    These are typed first so the shared-string table appends them before
#    the "Products\v2\Image" strings below (matches the author's order).
# ---------------------------------------------------------------------------
$cardBack = @(
    "v2\Card Back\Reactant.png",
    "v2\Card Back\Oxidant.png",
    "v2\Card Back\Reductant.png",
    "v2\Card Back\Acid.png",
    "v2\Card Back\Base.png",
    "v2\Card Back\Halogen.png",
    "v2\Card Back\Metals & Ylides.png",
    "v2\Card Back\Prot & Leav groups.png",
    "v2\Card Back\RXN COND & H2O.png",
    "v2\Card Back\Other.png",
    "v2\Card Back\Premium.png",
    "v2\Card Back\Product.png",
    "v2\Card Back\PDT.png"
)

$row = 200
foreach ($name in $cardBack) {
    $ws.Cells.Item($row, 5).Value = $name
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2) Remove the old running-total formulas in I175:I176 - the table grows
#    past them, so they get recreated further down (I188:I189).
# ---------------------------------------------------------------------------
$ws.Range("I175").ClearContents()
$ws.Range("I176").ClearContents()

# ---------------------------------------------------------------------------
# 3) New "Product" deck rows 174-186 (A/C/D/E/F/G), continuing the numbering
#    (18..30) and deck ("Product" / "Original Deck") from rows 157-173.
#    C holds the new "Products\v2\Image\N.png" shared strings.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 13; $i++) {
    $r = 174 + $i
    $ws.Cells.Item($r, 1).Value = 18 + $i
    $ws.Cells.Item($r, 3).Value = "v2\Products\v2\Image\" + ($i + 1) + ".png"
    $ws.Cells.Item($r, 4).Value = "Product"
    $ws.Cells.Item($r, 5).Value = "v2\Templates\Product.png"
    $ws.Cells.Item($r, 6).Value = "Original Deck"
    $ws.Cells.Item($r, 7).Value = 1
}

# ---------------------------------------------------------------------------
# 4) Recreate the running-total formulas below the extended table.
# ---------------------------------------------------------------------------
$ws.Range("I188").Formula = "=SUM(G2:G186)"
$ws.Range("I189").Formula = "=I188/9"

# ---------------------------------------------------------------------------
# 5) Column H got wider (content elsewhere made the best-fit column wider).
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 24.3

# ---------------------------------------------------------------------------
# 6) Leave the selection where the author ended up after the edit.
# ---------------------------------------------------------------------------
[void]$ws.Range("E213").Select()
